$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.798.72"
$ws.Range("E2").Value = "  +1.70%  "

$ws.Range("D3").Value = "1.663.25"
$ws.Range("E3").Value = "  +1.76%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'330.13"
$ws.Range("E5").Value = "  +8.30%  "

$ws.Range("E6").Value = "  -0.14%  "

$ws.Range("D7").Value = "'0.3647"
$ws.Range("E7").Value = "  +1.41%  "

$ws.Range("D8").Value = "'47.27"
$ws.Range("E8").Value = "  +1.01%  "

$ws.Range("D9").Value = "'0.3248"
$ws.Range("E9").Value = "  +0.47%  "

$ws.Range("D10").Value = "'1.137"
$ws.Range("E10").Value = "  +2.79%  "

$ws.Range("D11").Value = "'0.07064"
$ws.Range("E11").Value = "  +2.76%  "

$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("D13").Value = "'6.063"
$ws.Range("E13").Value = "  +2.57%  "

$ws.Range("D14").Value = "'19.51"
$ws.Range("E14").Value = "  +2.06%  "

$ws.Range("D15").Value = "1.665.45"
$ws.Range("E15").Value = "  +1.82%  "

$ws.Range("D16").Value = "'6.586"
$ws.Range("E16").Value = "  +1.09%  "

$ws.Range("D17").Value = "'0.00001048"
$ws.Range("E17").Value = "  +0.58%  "

$ws.Range("D18").Value = "'0.06637"
$ws.Range("E18").Value = "  +2.04%  "

$ws.Range("D20").Value = "'78.48"
$ws.Range("E20").Value = "  +2.69%  "

$ws.Range("D21").Value = "'5.924"
$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("D22").Value = "'15.78"
$ws.Range("E22").Value = "  +0.37%  "

$ws.Range("D23").Value = "'12.52"
$ws.Range("E23").Value = "  +4.54%  "

$ws.Range("D24").Value = "24.784.81"
$ws.Range("E24").Value = "  +1.87%  "

$ws.Range("D25").Value = "'2.465"
$ws.Range("E25").Value = "  +2.85%  "

$ws.Range("D26").Value = "'2.430"
$ws.Range("E26").Value = "  +4.68%  "

$ws.Range("D27").Value = "'148.45"
$ws.Range("E27").Value = "  +2.99%  "

$ws.Range("D28").Value = "'18.62"
$ws.Range("E28").Value = "  +0.39%  "

$ws.Range("D29").Value = "1.847.89"
$ws.Range("E29").Value = "  +1.68%  "

$ws.Range("D30").Value = "'125.76"
$ws.Range("E30").Value = "  +1.50%  "

$ws.Range("D31").Value = "'1.164"
$ws.Range("E31").Value = "  +3.45%  "

$ws.Range("D32").Value = "'4.068"
$ws.Range("E32").Value = "  +0.17%  "

$ws.Range("D33").Value = "'5.700"
$ws.Range("E33").Value = "  +1.20%  "

$ws.Range("D34").Value = "'0.08495"
$ws.Range("E34").Value = "  +1.86%  "

$ws.Range("D35").Value = "'1.636"
$ws.Range("E35").Value = "  -1.16%  "

$ws.Range("D36").Value = "'12.16"
$ws.Range("E36").Value = "  -1.05%  "

$ws.Range("D37").Value = "'5.158"
$ws.Range("E37").Value = "  +1.13%  "

$ws.Range("D38").Value = "'0.06164"
$ws.Range("E38").Value = "  +3.05%  "

$ws.Range("D39").Value = "'0.02282"
$ws.Range("E39").Value = "  +3.51%  "

$ws.Range("D40").Value = "'1.248"
$ws.Range("E40").Value = "  +4.17%  "

$ws.Range("D41").Value = "'0.2088"
$ws.Range("E41").Value = "  +3.01%  "

$ws.Range("D42").Value = "'8.223"
$ws.Range("E42").Value = "  +1.05%  "

$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("D44").Value = "'0.5919"
$ws.Range("E44").Value = "  +1.55%  "

$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "'3.850"
$ws.Range("E45").Value = "  +3.76%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'13.29"
$ws.Range("E46").Value = "  +6.28%  "

$ws.Range("D47").Value = "'0.5662"
$ws.Range("E47").Value = "  +2.54%  "

$ws.Range("D48").Value = "'125.73"
$ws.Range("E48").Value = "  +3.95%  "

$ws.Range("D49").Value = "'1.946"
$ws.Range("E49").Value = "  +1.51%  "

$ws.Range("D50").Value = "'0.06970"
$ws.Range("E50").Value = "  +1.32%  "

$ws.Range("E51").Value = "  +4.89%  "
